$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1, J1 - copy H1 style/format then set values
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I2:J69
$data = @{
    2 = @(8, 8)
    3 = @(6, 6)
    4 = @(7, 7)
    5 = @(8, 8)
    6 = @(7, 7)
    7 = @(8, 8)
    8 = @(6, 6)
    9 = @(7, 7)
    10 = @(10, 10)
    11 = @(9, 9)
    12 = @(9, 9)
    13 = @(9, 9)
    14 = @(7, 7)
    15 = @(6, 6)
    16 = @(8, 8)
    17 = @(8, 8)
    18 = @(7, 7)
    19 = @(7, 7)
    20 = @(6, 7)
    21 = @(7, 8)
    22 = @(9, 9)
    23 = @(7, 7)
    24 = @(6, 7)
    25 = @(7, 7)
    26 = @(9, 9)
    27 = @(8, 8)
    28 = @(8, 8)
    29 = @(8, 8)
    30 = @(7, 7)
    31 = @(8, 8)
    32 = @(7, 7)
    33 = @(9, 9)
    34 = @(8, 8)
    35 = @(8, 8)
    36 = @(7, 7)
    37 = @(6, 7)
    38 = @(9, 9)
    39 = @(8, 8)
    40 = @(7, 7)
    41 = @(6, 7)
    42 = @(9, 9)
    43 = @(8, 8)
    44 = @(7, 7)
    45 = @(9, 9)
    46 = @(2, 2)
    47 = @(8, 8)
    48 = @(7, 7)
    49 = @(7, 8)
    50 = @(7, 7)
    51 = @(7, 7)
    52 = @(8, 8)
    53 = @(8, 8)
    54 = @(8, 8)
    55 = @(6, 6)
    56 = @(5, 6)
    57 = @(7, 7)
    58 = @(6, 7)
    59 = @(6, 6)
    60 = @(7, 7)
    61 = @(5, 6)
    62 = @(6, 6)
    63 = @(7, 7)
    64 = @(9, 9)
    65 = @(6, 6)
    66 = @(5, 5)
    67 = @(7, 7)
    68 = @(8, 8)
    69 = @(9, 9)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 9).Value = $vals[0]
    $ws.Cells.Item($r, 10).Value = $vals[1]
}

Write-Host "Done"
